$d = $word.ActiveDocument

# --- Edit 1: mark the "Enter into the folder..." screenshot run as NoProof ---
# Locate the InlineShape matching the picture inserted after "Enter into the folder
# shown in below screenshot" (unique size 5943600 x 2221230 EMU == 468 x 174.9 pt).
$targetShape = $null
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shp = $d.InlineShapes.Item($i)
    if ([math]::Round($shp.Width, 1) -eq 468 -and [math]::Round($shp.Height, 1) -eq 174.9) {
        $targetShape = $shp
    }
}
if ($targetShape -ne $null) {
    $targetShape.Range.NoProofing = $true
}

# --- Edit 2: append the new "Enabling and Disabling..." section at the end of the
# document body, right after "To run crontab need to log into the raspberry pi over ssh " ---
$insertPoint = $d.Range($d.Content.End, $d.Content.End)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="IntenseQuote"/>
  </w:pPr>
  <w:r>
    <w:t>Enabling and Disabling Raspberry OS Feature to save energy</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t>2. USB AND ETHERNET</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>To shut off power for USB ports and Ethernet type the following into the</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>Raspberry Pi Terminal and press enter.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>echo '1-1' | sudo tee /sys/bus/usb/drivers/usb/unbind</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>To turn the power back on for USB ports and Ethernet type the following into the</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>Raspberry Pi Terminal and press enter.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>echo '1-1' | sudo tee /sys/bus/usb/drivers/usb/bind</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>3. WIFI AND BLUETOOTH</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>Disabling WiFi not only lower the overall power usage but also improves security</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>for your Raspberry Pi. There are a number of ways to do this. The easiest way is</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>to utilise the program | rfkill |. Type the following into the terminal and press enter</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t>for each line. Then on reboot it will be disabled.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>sudo apt install rfkill</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>sudo rfkill block wifi</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>If you want to re-enable WiFi simply type the following into the terminal.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>sudo rfkill unblock wifi</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>Disable Bluetooth</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>Internal</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>sudo rfkill block bluetooth</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>If you want to re-enable Bluetooth simply type the following into the terminal.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>sudo rfkill unblock Bluetooth</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t>4. DISABLE HDMI</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>Whenever you are remotely accessing the Raspberry Pi using SSH or running it</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>headless then likely there will be no need to run anything through the onboard</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>HDMI ports. This means you can disable power to these ports which will give a</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>small improvement to battery life (~20mA). Type the following into the terminal to</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>disable power to the HDMI ports. As soon as you press enter it will be disabled.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>sudo /opt/vc/bin/tvservice -o</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>If you want to re-enable HDMI connection simply type the following into the</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>terminal.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
  </w:pPr>
  <w:r>
    <w:t>sudo /opt/vc/bin/tvservice -p</w:t>
  </w:r>
</w:p>
'@
$insertPoint.InsertXML($xml)

Write-Output "edit complete"
